# The sheet holds one daily price record per row (rows 2..124). A new
# record for 2023-04-27 (serial 45043) needs to be inserted as the second
# record (row 5), pushing every subsequent record down by one row — the
# record that previously lived in the last row (124) ends up in the new
# last row (125).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 5; Excel shifts rows 5..124 down to
# 6..125 and grows the used range to A1:T125 automatically.
$ws.Rows("5:5").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C5").Value = "Ñuble"
$ws.Range("D5").Value = 45043
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100108
$ws.Range("H5").Value = "Tropicales y subtropicales"
$ws.Range("I5").Value = 100108002
$ws.Range("J5").Value = "Mango"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 70
$ws.Range("N5").Value = 7000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 7571
$ws.Range("Q5").Value = "$/bandeja 4 kilos"
$ws.Range("R5").Value = "Perú"
$ws.Range("S5").Value = 1893
$ws.Range("T5").Value = 4
